# Generate Report for Archive
#
# 1. Status flips from "Ready for handoff" to "In Translation". The
#    Overview sheet mirrors each language's status in columns E (zh-cn) and
#    F (de-de); each language sheet carries its own Status in column C.
#    (Cells are addressed directly rather than scanned via UsedRange.Cells,
#    which mutates unrelated cells that happen to share a value when
#    written to mid-iteration.)
# 2. The now-shorter status text lets the Status columns narrow.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the status columns. The COM layer snaps ColumnWidth to a 1/6
# character grid, so 12.5 is the closest achievable width to the recorded
# 13.4101845877511 (stored width = ColumnWidth + 5/6).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
